$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, shifting all existing rows down by one.
$ws.Rows.Item(1).Insert()

# Populate the new header row (unstyled, matching target formatting).
$ws.Range("A1").Value = "Q"
$ws.Range("B1").Value = "A1"
$ws.Range("C1").Value = "A2"
$ws.Range("D1").Value = "A3"
$ws.Range("E1").Value = "A4"
$ws.Range("F1").Value = "A"

$newQ = @{
    42 = "全面推进惩治和预防腐败体系建设，做到干部清正、政府清廉、()，反映了党推进反腐倡廉建设的坚定决心，对于提高党的建设科学化水平，凝聚党心民心，具有重要的指导意义。"
    43 = "党章规定的各项纪律都必须严格遵守和执行，而最首要、最核心的就是要严格遵守和执行党的( )"
    44 = "新时期党的干部工作的重要指导方针是( )"
    45 = "( )是马克思主义政党的基本政治观点和根本工作路线。"
    46 = "中国共产党是中国工人阶级的先锋队，同时是( )，是中国特色社会主义事业的领导核心，代表中国先进生产力的发展要求，代表中国先进文化的前进方向，代表中国最广大人民的根本利益。"
    47 = "党的最高理想和最终目标是( )。"
    48 = "党的根本宗旨是( )。"
    49 = "马克思列宁主义揭示了( )，它的基本原理是正确的，具有强大的生命力。"
    50 = "我国正处于并将长期处于( )。这是在经济文化落后的中国建设社会主义现代化不可逾越的历史阶段，需要上百年的时间。"
    51 = "中国共产党在社会主义初级阶段的( )是：领导和团结全国各族人民，以经济建设为中心，坚持四项基本原则，坚持改革开放，自力更生，艰苦创业，为把我国建设成为富强、民主、文明的社会主义现代化国家而奋斗。"
    52 = "党的最高领导机关，是( )和它所产生的中央委员会。"
    53 = "党的领导主要是( )的领导。"
    54 = "党的思想路线是一切从实际出发，理论联系实际，( )，在实践中检验真理和发展真理。"
    55 = "民主集中制是民主基础上的集中和集中指导下的民主相结合。它既是党的( )，也是群众路线在党的生活中的运用。"
    56 = "我们党的最大政治优势是( )，党执政后的最大危险是脱离群众。"
    57 = "在现阶段，我国社会的( )是人民日益增长的物质文化需要同落后的社会生产之间的矛盾。"
    58 = "必须尊重劳动、尊重知识、尊重人才、( )，做到发展为了人民、发展依靠人民、发展成果由人民共享。"
    59 = "党员必须自觉遵守党的纪律，模范遵守( )，严格保守党和国家的秘密，执行党的决定，服从组织分配，积极完成党的任务。"
    60 = "预备党员的( )同正式党员一样。"
    61 = "党员如果没有正当理由，连续( )不参加党的组织生活，或不交纳党费，或不做党所分配的工作，就被认为是自行脱党。"
}

foreach ($row in $newQ.Keys) {
    if ($row -ne 42) {
        $ws.Range("B$row").Copy() | Out-Null
        $ws.Range("A$row").PasteSpecial(-4122) | Out-Null
    }
    $ws.Range("A$row").Value = $newQ[$row]
}

$excel.CutCopyMode = 0

# Restore the view state described by the target workbook.
$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("G2").Select()
